$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 3.2
$ws.Range("S2").Value = 2.8
$ws.Range("AA3").Value = 290
$ws.Range("AB3").Value = 9.199999999999999
$ws.Range("AC3").Value = 12
$ws.Range("AD3").Value = 34
$ws.Range("AE3").Value = 150
$ws.Range("AF3").Value = 11
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 30
$ws.Range("AI3").Value = 130
$ws.Range("AJ3").Value = 18
$ws.Range("AK3").Value = 22
$ws.Range("AL3").Value = 48
$ws.Range("AM3").Value = 180
$ws.Range("AN3").Value = 11.5
$ws.Range("AO3").Value = 200
$ws.Range("G3").Value = 1.63
$ws.Range("M3").Value = 1.06
$ws.Range("P3").Value = 1.9
$ws.Range("Q3").Value = 1.91
$ws.Range("T3").Value = 1.96
$ws.Range("U3").Value = 1.81
$ws.Range("W3").Value = 2.58
$ws.Range("X3").Value = 18.5
$ws.Range("Y3").Value = 27
$ws.Range("Z3").Value = 75
$ws.Range("AL4").Value = 980
$ws.Range("W4").Value = 1.9
$ws.Range("F5").Value = 3.85
$ws.Range("G5").Value = 5.2
$ws.Range("I5").Value = 2.1
$ws.Range("J5").Value = 3.15
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 1.97
$ws.Range("T5").Value = 1.72
$ws.Range("U5").Value = 2.08
$ws.Range("V5").Value = 1.92
$ws.Range("W5").Value = 1.24
$ws.Range("I6").Value = 1.45
$ws.Range("J6").Value = 5.1
$ws.Range("R6").Value = 1.52
$ws.Range("V6").Value = 3.2
$ws.Range("AN7").Value = 2.92
$ws.Range("AO7").Value = 400
$ws.Range("H7").Value = 15.5
$ws.Range("I7").Value = 22
$ws.Range("N7").Value = 8
$ws.Range("R7").Value = 1.96
$ws.Range("U7").Value = 1.79
$ws.Range("V7").Value = 1.04
$ws.Range("G8").Value = 1.36
$ws.Range("H8").Value = 8.800000000000001
$ws.Range("K8").Value = 7.4
$ws.Range("W8").Value = 3.75
$ws.Range("AK9").Value = 38
$ws.Range("G9").Value = 2.58
$ws.Range("H9").Value = 2.96
$ws.Range("I9").Value = 3.75
$ws.Range("P9").Value = 1.79
$ws.Range("Q9").Value = 2.1
$ws.Range("S9").Value = 3.85
$ws.Range("T9").Value = 1.01
$ws.Range("U9").Value = 1.01
$ws.Range("V9").Value = 1.36
$ws.Range("W9").Value = 1.63
$ws.Range("AJ10").Value = 75
$ws.Range("AM10").Value = 130
$ws.Range("F10").Value = 3.2
$ws.Range("G10").Value = 3.65
$ws.Range("H10").Value = 2.34
$ws.Range("I10").Value = 2.6
$ws.Range("K10").Value = 3.55
$ws.Range("N10").Value = 3.35
$ws.Range("P10").Value = 1.8
$ws.Range("T10").Value = 1.81
$ws.Range("U10").Value = 2.08
$ws.Range("V10").Value = 1.62
$ws.Range("W10").Value = 1.38
$ws.Range("F11").Value = 1.95
$ws.Range("H11").Value = 3.85
$ws.Range("I11").Value = 4.4
$ws.Range("K11").Value = 4.1
$ws.Range("O11").Value = 1.24
$ws.Range("Q11").Value = 1.76
$ws.Range("R11").Value = 1.35
$ws.Range("S11").Value = 1.77
$ws.Range("T11").Value = 1.48
$ws.Range("V11").Value = 1.31
$ws.Range("G12").Value = 2
$ws.Range("Q12").Value = 1.4
$ws.Range("R12").Value = 1.88
$ws.Range("S12").Value = 1.94
$ws.Range("U12").Value = 2.9
$ws.Range("W12").Value = 2
$ws.Range("F13").Value = 2.38
$ws.Range("H13").Value = 2.6
$ws.Range("I13").Value = 2.92
$ws.Range("L13").Value = 1.22
$ws.Range("N13").Value = 6.4
$ws.Range("P13").Value = 2.82
$ws.Range("R13").Value = 1.74
$ws.Range("S13").Value = 2.1
$ws.Range("T13").Value = 1.45
$ws.Range("U13").Value = 2.8
$ws.Range("L14").Value = 1.25
$ws.Range("M14").Value = 1.02
$ws.Range("O14").Value = 1.2
$ws.Range("S14").Value = 2.32
$ws.Range("T14").Value = 1.59
$ws.Range("V14").Value = 1.29
$ws.Range("F15").Value = 1.65
$ws.Range("G15").Value = 1.77
$ws.Range("H15").Value = 4.4
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 4.4
$ws.Range("K15").Value = 5.2
$ws.Range("L16").Value = 1.24
$ws.Range("AA17").Value = 200
$ws.Range("AE17").Value = 90
$ws.Range("AH17").Value = 23
$ws.Range("AN17").Value = 5.3
$ws.Range("F17").Value = 1.44
$ws.Range("G17").Value = 1.53
$ws.Range("H17").Value = 6.4
$ws.Range("I17").Value = 8.6
$ws.Range("Q17").Value = 1.42
$ws.Range("S17").Value = 2.06
$ws.Range("U17").Value = 2.28
$ws.Range("V17").Value = 1.13
$ws.Range("W17").Value = 2.88
$ws.Range("Y17").Value = 40
$ws.Range("Z17").Value = 75
$ws.Range("F19").Value = 1.66
$ws.Range("G19").Value = 1.93
$ws.Range("I19").Value = 8.6
$ws.Range("J19").Value = 3.35
$ws.Range("K19").Value = 6.4
$ws.Range("N19").Value = 1.62
$ws.Range("P19").Value = 1.62
$ws.Range("Q19").Value = 1.98
$ws.Range("S19").Value = 1.98
$ws.Range("AA20").Value = 80
$ws.Range("AB20").Value = 8
$ws.Range("AC20").Value = 7.4
$ws.Range("AD20").Value = 16.5
$ws.Range("AE20").Value = 60
$ws.Range("AF20").Value = 16
$ws.Range("AG20").Value = 13.5
$ws.Range("AH20").Value = 24
$ws.Range("AI20").Value = 100
$ws.Range("AJ20").Value = 44
$ws.Range("AK20").Value = 40
$ws.Range("AL20").Value = 70
$ws.Range("AM20").Value = 220
$ws.Range("AN20").Value = 44
$ws.Range("AO20").Value = 80
$ws.Range("F20").Value = 2.42
$ws.Range("G20").Value = 2.76
$ws.Range("H20").Value = 3.25
$ws.Range("I20").Value = 3.8
$ws.Range("J20").Value = 2.86
$ws.Range("K20").Value = 3.3
$ws.Range("L20").Value = 1.55
$ws.Range("M20").Value = 1.12
$ws.Range("N20").Value = 2.52
$ws.Range("O20").Value = 1.52
$ws.Range("P20").Value = 1.52
$ws.Range("Q20").Value = 2.54
$ws.Range("R20").Value = 1.18
$ws.Range("S20").Value = 4.7
$ws.Range("T20").Value = 2.06
$ws.Range("U20").Value = 1.76
$ws.Range("V20").Value = 1.35
$ws.Range("W20").Value = 1.56
$ws.Range("X20").Value = 8.800000000000001
$ws.Range("Y20").Value = 10
$ws.Range("Z20").Value = 24
$ws.Range("F21").Value = 1.89
$ws.Range("N21").Value = 3.35
$ws.Range("Q21").Value = 2.12
$ws.Range("V21").Value = 1.25
$ws.Range("O23").Value = 1.53
$ws.Range("S23").Value = 5.6
$ws.Range("J24").Value = 2.7
$ws.Range("AH25").Value = 19.5
$ws.Range("I25").Value = 3.25
$ws.Range("V25").Value = 1.44
$ws.Range("AM27").Value = 310
